$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.058.70"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "3.412.69"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.88"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.98"
$ws.Range("E6").Value = "  +5.10%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.477"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.54"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.127"
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.393"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Value = "3.991.16"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000179"
$ws.Range("E14").Value = "  +4.39%  "
$ws.Range("D15").Value = "3.413.54"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.71"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "62.097.91"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.23"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.57"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.84"
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.35"
$ws.Range("E21").Value = "  +6.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.569"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000131"
$ws.Range("E23").Value = "  +14.30%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.541.68"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.64"
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.77"
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.60"
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.34"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.162"
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.20"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.57"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.49"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.442.42"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.59"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.99"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.64"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0795"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  +13.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.28"
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.786"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.26"
$ws.Range("E45").Value = "  +6.68%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.47"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.82"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.95"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.51"
$ws.Range("E49").Value = "  +3.79%  "
$ws.Range("D50").Value = "2.364.09"
$ws.Range("E50").Value = "  +8.79%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.885"
$ws.Range("E51").Value = "  -1.41%  "
